$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23 (shifts existing rows 23..124 down to 24..125,
# carrying formatting from the surrounding rows, matching Excel's default
# "insert" behaviour).
$ws.Rows("23").Insert()

# Populate the newly inserted row 23 with the new price-record data.
$ws.Cells.Item(23, 1).Value = 11
$ws.Cells.Item(23, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(23, 3).Value = "Bíobío"
$ws.Cells.Item(23, 4).Value = 44859
$ws.Cells.Item(23, 5).Value = 8
$ws.Cells.Item(23, 6).Value = 100112001
$ws.Cells.Item(23, 7).Value = "Berenjena"
$ws.Cells.Item(23, 8).Value = "Sin especificar"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 220
$ws.Cells.Item(23, 11).Value = 10000
$ws.Cells.Item(23, 12).Value = 12000
$ws.Cells.Item(23, 13).Value = 10909
$ws.Cells.Item(23, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(23, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(23, 16).Value = 182
$ws.Cells.Item(23, 17).Value = 60
$ws.Cells.Item(23, 18).Value = "Hortaliza"
